# TAC-3104, Fix start trip date issue
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference No (A2): 3 -> 4
$ws.Range("A2").Value = 4

# Trip Pick up Date Start (B2): was text "03/19/2022", fix to an actual date
# value 3/20/2022 (date serial 44640) formatted as m/d/yyyy, left aligned.
$ws.Range("B2").Value = 44640
$ws.Range("B2").NumberFormat = "[$-1010000]m/d/yyyy;@"
$ws.Range("B2").HorizontalAlignment = -4131

# Move the active selection to B2 (was G2)
$ws.Range("B2").Select()
